$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.798.60'
$ws.Range('E2').Value = '  +0.48%  '
$ws.Range('D3').Value = '1.644.81'
$ws.Range('E3').Value = '  +0.07%  '
$ws.Range('E4').Value = '  +0.35%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '216.73'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.62%  '
$ws.Range('E6').Value = '  -0.55%  '
$ws.Range('E7').Value = '  +0.29%  '
$ws.Range('E8').Value = '  -0.21%  '
$ws.Range('E9').Value = '  +0.12%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.21'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.32%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0846'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.43%  '
$ws.Range('D12').Value = '1.630.19'
$ws.Range('E12').Value = '  -1.01%  '
$ws.Range('E13').Value = '  -0.72%  '
$ws.Range('E14').Value = '  -0.22%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '64.70'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.79%  '
$ws.Range('D16').Value = '26.811.06'
$ws.Range('E16').Value = '  +0.42%  '
$ws.Range('E17').Value = '  -1.32%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '214.38'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.94%  '
$ws.Range('E19').Value = '  +0.29%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.40'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.13%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '2.44'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +7.76%  '
$ws.Range('E22').Value = '  -0.24%  '
$ws.Range('E23').Value = '  -1.82%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '146.22'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.03%  '
$ws.Range('E25').Value = '  +0.18%  '
$ws.Range('E26').Value = '  -1.40%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.20'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.28%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '15.69'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.27%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.0509'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.66%  '
$ws.Range('E30').Value = '  +0.46%  '
$ws.Range('E31').Value = '  -1.04%  '
$ws.Range('E32').Value = '  -1.34%  '
$ws.Range('D33').Value = '1.286.98'
$ws.Range('E33').Value = '  +1.02%  '
$ws.Range('E34').Value = '  -0.27%  '
$ws.Range('E35').Value = '  +1.40%  '
$ws.Range('E36').Value = '  -1.00%  '
$ws.Range('E37').Value = '  +0.57%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.822'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.12%  '
$ws.Range('E39').Value = '  +0.30%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.807'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.32%  '
$ws.Range('E41').Value = '  -0.74%  '
$ws.Range('E42').Value = '  -2.11%  '
$ws.Range('D43').Value = '1.784.01'
$ws.Range('E43').Value = '  +0.03%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '61.44'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +3.01%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '91.89'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.27%  '
$ws.Range('E46').Value = '  +0.43%  '
$ws.Range('B47').Value = 'Cronos'
$ws.Range('C47').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0520'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.76%  '
$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '7.66'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -2.31%  '
$ws.Range('B49').Value = 'Algorand'
$ws.Range('C49').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0969'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.08%  '
$ws.Range('B50').Value = 'Mantle'
$ws.Range('C50').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.407'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.10%  '
$ws.Range('B51').Value = 'USDD'
$ws.Range('C51').Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.01'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.34%  '
